# "edit the dataset 1" - update the Location column (D) on sheet "s1" with
# corrected / standardized country names, matching a repeating 10-value
# cycle: Sweden, United Kingdom, Germany, France, Italy, Spain,
# "united states of america ", Russia, Denmark, Finland.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$locations = @(
  "Sweden",
  "United Kingdom",
  "Germany",
  "France",
  "Italy",
  "Spain",
  "united states of america ",
  "Russia",
  "Denmark",
  "Finland"
)

for ($row = 2; $row -le 30; $row++) {
  $value = $locations[($row - 2) % $locations.Length]
  $ws.Cells.Item($row, 4).Value = $value
}

# The new longest entry ("united states of america ") makes column D much
# wider than before; set its width accordingly (was bestFit to ~11.29,
# now bestFit to ~21.86 characters).
$ws.Columns.Item(4).ColumnWidth = 21

# Move/restore the active selection like the author left it.
$ws.Range("D40").Select() | Out-Null
